# Weekly update: insert a new price record as row 36, pushing the
# existing rows 36-58 down to 37-59 (last row duplicated data shifts too).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 36.
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with this week's data.
$ws.Range("A36").Value = 4
$ws.Range("B36").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C36").Value = "Los Lagos"
$ws.Range("D36").Value = 44875
$ws.Range("E36").Value = 10
$ws.Range("F36").Value = 300000000
$ws.Range("G36").Value = "Espárragos"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 200
$ws.Range("K36").Value = 1500
$ws.Range("L36").Value = 1700
$ws.Range("M36").Value = 1600
$ws.Range("N36").Value = "`$/kilo"
$ws.Range("O36").Value = "Provincia de Linares"
$ws.Range("P36").Value = 1600
$ws.Range("Q36").Value = 1
$ws.Range("R36").Value = "Hortaliza"
